$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: cache every "numeric-looking" text value that we'll need again
# later, copying it (as VALUES, so it keeps its shared-string / text nature)
# from wherever it currently lives in the sheet into a scratch area far away
# (row 200+) before any of those source cells get overwritten.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy(); $ws.Range("A200").PasteSpecial(-4163)   # "400"
$ws.Range("B3").Copy(); $ws.Range("A201").PasteSpecial(-4163)   # "402"
$ws.Range("C2").Copy(); $ws.Range("A202").PasteSpecial(-4163)   # "01"
$ws.Range("D2").Copy(); $ws.Range("A203").PasteSpecial(-4163)   # "17502"
$ws.Range("D3").Copy(); $ws.Range("A204").PasteSpecial(-4163)   # "8130"
$ws.Range("E2").Copy(); $ws.Range("A205").PasteSpecial(-4163)   # "2"
$ws.Range("E3").Copy(); $ws.Range("A206").PasteSpecial(-4163)   # "4"
$ws.Range("J2").Copy(); $ws.Range("A207").PasteSpecial(-4163)   # "121"

# "0700" never existed as its own text value before (it was embedded inside
# "0700-0900pm"). Derive it with a formula, then bake the formula result into
# a plain text value via a self copy/paste-special so it becomes a normal
# shared-string cell instead of a cached formula.
$ws.Range("A208").Formula = "=LEFT(H2,4)"
$ws.Range("A208").Copy(); $ws.Range("A208").PasteSpecial(-4163)

$c400   = $ws.Range("A200")
$c402   = $ws.Range("A201")
$c01    = $ws.Range("A202")
$c17502 = $ws.Range("A203")
$c8130  = $ws.Range("A204")
$c2     = $ws.Range("A205")
$c4     = $ws.Range("A206")
$c121   = $ws.Range("A207")
$c0700  = $ws.Range("A208")

# ---------------------------------------------------------------------------
# Phase 2: numeric header row (row 1, columns B..M = 0..11) and numeric
# index column (A2:A5 = 0..3). B1:K1 and A2:A3 already carry the bordered /
# bold "header" style from the source file; L1, M1, A4 and A5 are brand new
# cells so copy that formatting across (format-only paste, so no new style
# entries are created) before dropping the numbers in.
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11

$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# Phase 3: the actual data grid (text columns). Plain words/names are typed
# directly; anything that looks like a number is pasted in as a value from
# the scratch cache above so Excel keeps treating it as text.
# ---------------------------------------------------------------------------

# Row 2 -- B2/C2/D2/E2/F2 already hold the correct text ("400","01","17502","2","F")
$ws.Range("G2").Value = "W"
$ws.Range("H2").Value = "Lat Am Studies Capstone"
$c0700.Copy();  $ws.Range("I2").PasteSpecial(-4163)
$ws.Range("J2").Value = "0900pm"
$ws.Range("K2").Value = "KING"
$c121.Copy();   $ws.Range("L2").PasteSpecial(-4163)
$ws.Range("M2").Value = "Mani Kristina"

# Row 3
$c400.Copy();   $ws.Range("B3").PasteSpecial(-4163)
$ws.Range("C3").Value = "W"
$c0700.Copy();  $ws.Range("E3").PasteSpecial(-4163)
$ws.Range("F3").Value = "0900pm"
$ws.Range("G3").Value = "KING"
$c121.Copy();   $ws.Range("H3").PasteSpecial(-4163)
$ws.Range("I3").Value = "O'Connor Patrick"

# Row 4
$c402.Copy();   $ws.Range("B4").PasteSpecial(-4163)
$c01.Copy();    $ws.Range("C4").PasteSpecial(-4163)
$c8130.Copy();  $ws.Range("D4").PasteSpecial(-4163)
$c4.Copy();     $ws.Range("E4").PasteSpecial(-4163)
$ws.Range("F4").Value = "F"
$ws.Range("H4").Value = "Senior"
$ws.Range("I4").Value = "Honors"
$ws.Range("J4").Value = "TBA"
$ws.Range("K4").Value = "TBA"
$ws.Range("L4").Value = "Mani Kristina"

# Row 5
$c402.Copy();   $ws.Range("B5").PasteSpecial(-4163)
$ws.Range("C5").Value = "TBA"
$ws.Range("D5").Value = "TBA"
$ws.Range("E5").Value = "O'Connor Patrick"

# ---------------------------------------------------------------------------
# Phase 4: drop the scratch helper cells so they don't show up in the saved
# sheet / used range.
# ---------------------------------------------------------------------------
$ws.Range("A200:A208").Clear()
